$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header A1 from "Command" to "System"
$ws.Range("A1").Value = "System"

# Add note to C2
$ws.Range("C2").Value = "Need Restructure"

# Append new row 26
$ws.Range("A26").Value = "Prefix"
$ws.Range("B26").Value = "Event"
$ws.Range("C26").Value = "s!"
$ws.Range("F26").Value = "Complete 0.3.0.1b"
